$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Day 14" row (row 15) ---
# New shared strings must be introduced in the same left-to-right order
# they appear in the target file: A15, C15, D15, E15, F15, then the "S "
# status variant used by G4/G6.

# A15: day label
$ws.Range("A15").Value = "Day 14"

# C15 / D15 / E15: hyperlinked question files for Day 14.
# Hyperlinks.Add's TextToDisplay both seeds the cell text and the cached
# <hyperlink display="…"> attribute — Excel actually caches display="" as
# the target Address rather than the friendly file name, so add the link
# with TextToDisplay = Address first, then overwrite the cell text with the
# friendly file name afterwards (this does not disturb the already-written
# display="" attribute, it only rewrites the shared string in place).
$dup = "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 14/Contains Duplicate II.py"
$ws.Hyperlinks.Add($ws.Range("C15"), $dup, "", "Contains Duplicate II.py", $dup) | Out-Null
$ws.Range("C15").Value = "Contains Duplicate II.py"
$ws.Range("C15").Style = "Hyperlink"

$temps = "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 14/Daily Temperatures.py"
$ws.Hyperlinks.Add($ws.Range("D15"), $temps, "", "Daily Temperatures.py", $temps) | Out-Null
$ws.Range("D15").Value = "Daily Temperatures.py"
$ws.Range("D15").Style = "Hyperlink"

$rpn = "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 14/Evaluate Reverse Polish Notation.py"
$ws.Hyperlinks.Add($ws.Range("E15"), $rpn, "", "Evaluate Reverse Polish Notation.py", $rpn) | Out-Null
$ws.Range("E15").Value = "Evaluate Reverse Polish Notation.py"
$ws.Range("E15").Style = "Hyperlink"

# F15 / G15: plain text cells ("S" already exists in the shared string table)
$ws.Range("F15").Value = "Stack, HashMap"
$ws.Range("G15").Value = "S"

# B15: date — copy the date format from B14 (keeps the existing numFmt/style
# instead of minting a new one), then write the raw serial value for
# 2025-06-08 (one day after Day 13's 2025-06-07). No new shared string.
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("B15").Value = 45816

# --- Day 4 / Day 6 status cells now use the "S " (trailing space) variant ---
$ws.Range("G4").Value = "S "
$ws.Range("G6").Value = "S "

# --- Dimension grows to A1:H15 automatically; move the active selection ---
$ws.Range("D10").Select()
